$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Fix the Avg_Time_ms values for the first two benchmark rows
# (Rows = 5000 and Rows = 10000) in the InsertBinarySort timing data.
$ws.Range("D2").Value = 20.578025
$ws.Range("D3").Value = 83.29548000000001

